$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G403").Value = 999329.6854957336
$ws.Range("G404").Value = 2003892.370953069
$ws.Range("G405").Value = 3013576.934744317
$ws.Range("G406").Value = 4028275.217237662
$ws.Range("G407").Value = 5047882.39340258
$ws.Range("G408").Value = 6072297.344322637
$ws.Range("G409").Value = 7101423.026479859
$ws.Range("G410").Value = 8135166.837623404
$ws.Range("G411").Value = 9173440.977989187
$ws.Range("G412").Value = 10216162.80557838
$ws.Range("G413").Value = 10216155.95661702
$ws.Range("H413").Value = 997383.0446407206
$ws.Range("G414").Value = 10216147.45734638
$ws.Range("H414").Value = 1998862.792524885
$ws.Range("G415").Value = 10216136.93408722
$ws.Range("H415").Value = 3004377.574296403
$ws.Range("G416").Value = 10216123.93447644
$ws.Range("H416").Value = 4013871.978856502
$ws.Range("G417").Value = 10216107.91229236
$ws.Range("H417").Value = 5027297.157387884
$ws.Range("G418").Value = 10216088.2096637
$ws.Range("H418").Value = 6044611.115497725
$ws.Range("G419").Value = 10216064.03627614
$ws.Range("H419").Value = 7065778.99212087
$ws.Range("G420").Value = 10216034.44514672
$ws.Range("H420").Value = 8090773.3237717
$ws.Range("G421").Value = 10215998.30448963
$ws.Range("H421").Value = 9119574.292677579
$ws.Range("G422").Value = 10215954.26514854
$ws.Range("H422").Value = 10152169.95726554
$ws.Range("G423").Value = 10215900.72302018
$ws.Range("H423").Value = 11188556.46341213
$ws.Range("G424").Value = 10215835.77584281
$ws.Range("H424").Value = 12228738.23479548
$ws.Range("G425").Value = 10215757.17367182
$ws.Range("H425").Value = 13272728.1406202
$ws.Range("G426").Value = 10215662.26231411
$ws.Range("H426").Value = 14320547.63890641
$ws.Range("G427").Value = 10215547.91894458
$ws.Range("H427").Value = 15372226.89345749
$ws.Range("G428").Value = 10215410.47908301
$ws.Range("H428").Value = 16427804.862538
$ws.Range("G429").Value = 10215245.6540698
$ws.Range("H429").Value = 17487329.35720627
$ws.Range("G430").Value = 10215048.43814639
$ws.Range("H430").Value = 18550857.06716072
$ws.Range("G431").Value = 10214813.00422235
$ws.Range("H431").Value = 19618453.55186832
$ws.Range("G432").Value = 10214532.58739918
$ws.Range("H432").Value = 20690193.19465396
$ws.Range("G433").Value = 10214199.35532242
$ws.Range("H433").Value = 21766159.1173418
$ws.Range("G434").Value = 10213804.26445258
$ws.Range("H434").Value = 22846443.05295034
$ws.Range("G435").Value = 10213336.90138327
$ws.Range("H435").Value = 23931145.17386251
$ws.Range("G436").Value = 10212785.30839615
$ws.Range("H436").Value = 25020373.87281014
$ws.Range("G437").Value = 10212135.79252854
$ws.Range("H437").Value = 26114245.49394601
$ws.Range("G438").Value = 10211372.71754512
$ws.Range("H438").Value = 27212884.01121197
$ws.Range("G439").Value = 10210478.27835205
$ws.Range("H439").Value = 28316420.6511682
$ws.Range("G440").Value = 10209432.25757374
$ws.Range("H440").Value = 29424993.45741284
$ws.Range("G441").Value = 10208211.76423088
$ws.Range("H441").Value = 30538746.79371208
$ws.Range("G442").Value = 10206790.9547168
$ws.Range("H442").Value = 31657830.78297016
$ws.Range("G443").Value = 10205140.73656759
$ws.Range("H443").Value = 32782400.67920551
$ws.Range("G444").Value = 10203228.45586285
$ws.Range("H444").Value = 33912616.16977102
$ws.Range("G445").Value = 10201017.56947617
$ws.Range("H445").Value = 35048640.60516031
$ws.Range("G446").Value = 10198467.30381959
$ws.Range("H446").Value = 36190640.15388992
$ws.Range("G447").Value = 10195532.30218984
$ws.Range("H447").Value = 37338782.88013972
$ws.Range("G448").Value = 10192162.26332577
$ws.Range("H448").Value = 38493237.74207944
$ws.Range("G449").Value = 10188301.57431997
$ws.Range("H449").Value = 39654173.50911298
$ws.Range("G450").Value = 10183888.94158877
$ws.Range("H450").Value = 40821757.59663185
$ws.Range("G451").Value = 10178857.0241861
$ws.Range("H451").Value = 41996154.81730614
$ws.Range("G452").Value = 10173132.07433967
$ws.Range("H452").Value = 43177526.04844242
